# Auto-generated Excel COM-interop script
# Applies numeric updates to leve profit calculation columns (H-N)
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33 (Leve Item ID 5512)
$ws.Range("H33").Value = 302.35715
$ws.Range("I33").Value = 302.35715
$ws.Range("K33").Value = 302.35715
$ws.Range("M33").Value = -73.35714999999999

# Row 38 (Leve Item ID 4599)
$ws.Range("H38").Value = 109.72727
$ws.Range("I38").Value = 50.7
$ws.Range("J38").Value = 700
$ws.Range("K38").Value = 152.1
$ws.Range("L38").Value = 2100
$ws.Range("M38").Value = 219.9
$ws.Range("N38").Value = -2844

# Row 58 (Leve Item ID 4606)
$ws.Range("H58").Value = 2184.1667
$ws.Range("I58").Value = 621
$ws.Range("J58").Value = 10000
$ws.Range("K58").Value = 1863
$ws.Range("L58").Value = 30000
$ws.Range("M58").Value = -1713
$ws.Range("N58").Value = -30300

# Row 125 (Leve Item ID 36228)
$ws.Range("H125").Value = 623.1
$ws.Range("I125").Value = 438.5
$ws.Range("K125").Value = 3946.5
$ws.Range("M125").Value = -1486.5

# Row 129 (Leve Item ID 36115)
$ws.Range("H129").Value = 295377.7
$ws.Range("J129").Value = 304304.38
$ws.Range("L129").Value = 912913.14
$ws.Range("N129").Value = -922913.14

# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 2851.8125
$ws.Range("I132").Value = 2927.6553
$ws.Range("K132").Value = 8782.965899999999
$ws.Range("M132").Value = -6252.965899999999

# Row 137 (Leve Item ID 44013)
$ws.Range("H137").Value = 1616.35
$ws.Range("I137").Value = 1229.7273
$ws.Range("K137").Value = 3689.1819
$ws.Range("M137").Value = -1139.1819

# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 12502733
$ws.Range("J138").Value = 3154.7795
$ws.Range("L138").Value = 9464.3385
$ws.Range("N138").Value = -19744.3385

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 4024.2134
$ws.Range("I32").Value = 3839.4177
$ws.Range("J32").Value = 5484.1
$ws.Range("K32").Value = 3839.4177
$ws.Range("L32").Value = 5484.1
$ws.Range("M32").Value = -3552.4177
$ws.Range("N32").Value = -6058.1

# Row 41 (Leve Item ID 2501)
$ws.Range("H41").Value = 3033.6
$ws.Range("I41").Value = 3033.6
$ws.Range("K41").Value = 3033.6
$ws.Range("M41").Value = -2619.6

# Row 61 (Leve Item ID 43999)
$ws.Range("H61").Value = 440061.44
$ws.Range("I61").Value = 515353
$ws.Range("J61").Value = 860.6667
$ws.Range("K61").Value = 515353
$ws.Range("L61").Value = 860.6667
$ws.Range("M61").Value = -515141
$ws.Range("N61").Value = -1284.6667

# Row 132 (Leve Item ID 43997)
$ws.Range("H132").Value = 11111.294
$ws.Range("I132").Value = 1279.7142
$ws.Range("K132").Value = 3839.1426
$ws.Range("M132").Value = -1309.1426

# Row 136 (Leve Item ID 43999)
$ws.Range("H136").Value = 440061.44
$ws.Range("I136").Value = 515353
$ws.Range("J136").Value = 860.6667
$ws.Range("K136").Value = 1546059
$ws.Range("L136").Value = 2582.0001
$ws.Range("M136").Value = -1543509
$ws.Range("N136").Value = -7682.0001

$ws = $wb.Worksheets.Item("BSM")
# Row 20 (Leve Item ID 14149)
$ws.Range("H20").Value = 1435.238
$ws.Range("I20").Value = 1605.8572
$ws.Range("J20").Value = 1094
$ws.Range("K20").Value = 1605.8572
$ws.Range("L20").Value = 1094
$ws.Range("M20").Value = -1358.8572
$ws.Range("N20").Value = -1588

# Row 76 (Leve Item ID 10630)
$ws.Range("H76").Value = 8000
$ws.Range("J76").Value = 8000
$ws.Range("L76").Value = 8000
$ws.Range("N76").Value = -8630

# Row 79 (Leve Item ID 10630)
$ws.Range("H79").Value = 8000
$ws.Range("J79").Value = 8000
$ws.Range("L79").Value = 8000
$ws.Range("N79").Value = -10184

# Row 134 (Leve Item ID 43998)
$ws.Range("H134").Value = 4008.0605
$ws.Range("I134").Value = 4405.846
$ws.Range("J134").Value = 2530.5715
$ws.Range("K134").Value = 13217.538
$ws.Range("L134").Value = 7591.7145
$ws.Range("M134").Value = -10682.538
$ws.Range("N134").Value = -12661.7145

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 2301.5967
$ws.Range("I31").Value = 1358.4131
$ws.Range("K31").Value = 1358.4131
$ws.Range("M31").Value = -1063.4131

# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 2301.5967
$ws.Range("I34").Value = 1358.4131
$ws.Range("K34").Value = 1358.4131
$ws.Range("M34").Value = -1156.4131

# Row 58 (Leve Item ID 44021)
$ws.Range("H58").Value = 26273.75
$ws.Range("I58").Value = 1362.4117
$ws.Range("K58").Value = 1362.4117
$ws.Range("M58").Value = -1159.4117

# Row 105 (Leve Item ID 19928)
$ws.Range("H105").Value = 10417589
$ws.Range("I105").Value = 13889641
$ws.Range("K105").Value = 13889641
$ws.Range("M105").Value = -13887894

# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 1870.2727
$ws.Range("I132").Value = 1425.0698
$ws.Range("K132").Value = 4275.2094
$ws.Range("M132").Value = -1745.2094

# Row 134 (Leve Item ID 44020)
$ws.Range("H134").Value = 906.36365
$ws.Range("I134").Value = 752.1579
$ws.Range("K134").Value = 2256.4737
$ws.Range("M134").Value = 278.5263

# Row 136 (Leve Item ID 44021)
$ws.Range("H136").Value = 26273.75
$ws.Range("I136").Value = 1362.4117
$ws.Range("K136").Value = 4087.2351
$ws.Range("M136").Value = -1537.2351

$ws = $wb.Worksheets.Item("CUL")
# Row 18 (Leve Item ID 36056)
$ws.Range("H18").Value = 198.77777
$ws.Range("I18").Value = 198.77777
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 596.33331
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = -427.33331
$ws.Range("N18").ClearContents()

# Row 126 (Leve Item ID 36045)
$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 15000
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -24880

# Row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 695.89
$ws.Range("J131").Value = 712.9239
$ws.Range("L131").Value = 2138.7717
$ws.Range("N131").Value = -12218.7717

$ws = $wb.Worksheets.Item("GSM")
# Row 11 (Leve Item ID 4422)
$ws.Range("H11").Value = 5202000.5
$ws.Range("I11").Value = 5500000
$ws.Range("J11").Value = 4010002
$ws.Range("K11").Value = 5500000
$ws.Range("L11").Value = 4010002
$ws.Range("M11").Value = -5499861
$ws.Range("N11").Value = -4010280

# Row 70 (Leve Item ID 14146)
$ws.Range("H70").Value = 10733.444
$ws.Range("I70").Value = 10246.714
$ws.Range("K70").Value = 10246.714
$ws.Range("M70").Value = -9976.714

# Row 73 (Leve Item ID 14146)
$ws.Range("H73").Value = 10733.444
$ws.Range("I73").Value = 10246.714
$ws.Range("K73").Value = 10246.714
$ws.Range("M73").Value = -9310.714

# Row 126 (Leve Item ID 36184)
$ws.Range("H126").Value = 5156.56
$ws.Range("I126").Value = 4243.75
$ws.Range("J126").Value = 6779.3335
$ws.Range("K126").Value = 12731.25
$ws.Range("L126").Value = 20338.0005
$ws.Range("M126").Value = -10261.25
$ws.Range("N126").Value = -25278.0005

# Row 132 (Leve Item ID 44008)
$ws.Range("H132").Value = 25519.39
$ws.Range("I132").Value = 4068.8572
$ws.Range("K132").Value = 12206.5716
$ws.Range("M132").Value = -9676.571599999999

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Leve Item ID 36249)
$ws.Range("H7").Value = 50005480
$ws.Range("I7").Value = 71431560
$ws.Range("J7").Value = 11284.167
$ws.Range("K7").Value = 71431560
$ws.Range("L7").Value = 11284.167
$ws.Range("M7").Value = -71431448
$ws.Range("N7").Value = -11508.167

# Row 55 (Leve Item ID 5284)
$ws.Range("H55").Value = 185.5
$ws.Range("I55").Value = 185.38461
$ws.Range("J55").Value = 185.8
$ws.Range("K55").Value = 185.38461
$ws.Range("L55").Value = 185.8
$ws.Range("M55").Value = -12.38461000000001
$ws.Range("N55").Value = -531.8

# Row 61 (Leve Item ID 27740)
$ws.Range("H61").Value = 7500
$ws.Range("I61").Value = 3366.6667
$ws.Range("J61").Value = 9566.666999999999
$ws.Range("K61").Value = 3366.6667
$ws.Range("L61").Value = 9566.666999999999
$ws.Range("M61").Value = -3164.6667
$ws.Range("N61").Value = -9970.666999999999

# Row 93 (Leve Item ID 19993)
$ws.Range("H93").Value = 1817.1666
$ws.Range("I93").Value = 1680.6
$ws.Range("J93").Value = 2500
$ws.Range("K93").Value = 1680.6
$ws.Range("L93").Value = 2500
$ws.Range("M93").Value = -432.5999999999999
$ws.Range("N93").Value = -4996

# Row 113 (Leve Item ID 27740)
$ws.Range("H113").Value = 7500
$ws.Range("I113").Value = 3366.6667
$ws.Range("J113").Value = 9566.666999999999
$ws.Range("K113").Value = 3366.6667
$ws.Range("L113").Value = 9566.666999999999
$ws.Range("M113").Value = -1196.6667
$ws.Range("N113").Value = -13906.667

# Row 126 (Leve Item ID 36249)
$ws.Range("H126").Value = 50005480
$ws.Range("I126").Value = 71431560
$ws.Range("J126").Value = 11284.167
$ws.Range("K126").Value = 214294680
$ws.Range("L126").Value = 33852.501
$ws.Range("M126").Value = -214292210
$ws.Range("N126").Value = -38792.501

# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 464707.03
$ws.Range("I132").Value = 483263.3
$ws.Range("K132").Value = 1449789.9
$ws.Range("M132").Value = -1447259.9

# Row 136 (Leve Item ID 44060)
$ws.Range("H136").Value = 999.6667
$ws.Range("I136").Value = 932.76666
$ws.Range("K136").Value = 2798.29998
$ws.Range("M136").Value = -248.2999799999998

$ws = $wb.Worksheets.Item("WVR")
# Row 62 (Leve Item ID 12589)
$ws.Range("H62").Value = 3690.9092
$ws.Range("I62").Value = 3175.25
$ws.Range("J62").Value = 3985.5715
$ws.Range("K62").Value = 3175.25
$ws.Range("L62").Value = 3985.5715
$ws.Range("M62").Value = -2551.25
$ws.Range("N62").Value = -5233.5715

# Row 64 (Leve Item ID 11036)
$ws.Range("H64").Value = 8888
$ws.Range("I64").Value = 8888
$ws.Range("K64").Value = 8888
$ws.Range("M64").Value = -8640

# Row 65 (Leve Item ID 12589)
$ws.Range("H65").Value = 3690.9092
$ws.Range("I65").Value = 3175.25
$ws.Range("J65").Value = 3985.5715
$ws.Range("K65").Value = 15876.25
$ws.Range("L65").Value = 19927.8575
$ws.Range("M65").Value = -12756.25
$ws.Range("N65").Value = -26167.8575

# Row 67 (Leve Item ID 11036)
$ws.Range("H67").Value = 8888
$ws.Range("I67").Value = 8888
$ws.Range("K67").Value = 8888
$ws.Range("M67").Value = -8030

# Row 75 (Leve Item ID 11957)
$ws.Range("H75").Value = 12000
$ws.Range("J75").Value = 12000
$ws.Range("L75").Value = 12000
$ws.Range("N75").Value = -13872

# Row 78 (Leve Item ID 11957)
$ws.Range("H78").Value = 12000
$ws.Range("J78").Value = 12000
$ws.Range("L78").Value = 36000
$ws.Range("N78").Value = -45360

# Row 126 (Leve Item ID 36210)
$ws.Range("H126").Value = 1226.0952
$ws.Range("I126").Value = 1162.4
$ws.Range("K126").Value = 3487.2
$ws.Range("M126").Value = -1017.2

# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 495.74545
$ws.Range("I132").Value = 469.1698
$ws.Range("K132").Value = 1407.5094
$ws.Range("M132").Value = 1122.4906

# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 18111874
$ws.Range("I136").Value = 24006904
$ws.Range("K136").Value = 72020712
$ws.Range("M136").Value = -72018162

Write-Host "Applied all leve profit updates."